# Update the multiplication problems in the table to the new values
# from the commit, addressing each cell directly by (row, column) so
# that duplicate expressions (e.g. "904×5=" appearing twice) are each
# replaced with their own distinct target value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="140×6="},
    @{Row=1;  Col=2; Text="345×6="},
    @{Row=1;  Col=3; Text="253×2="},
    @{Row=1;  Col=4; Text="394×2="},
    @{Row=1;  Col=5; Text="393×2="},

    @{Row=5;  Col=1; Text="895×7="},
    @{Row=5;  Col=2; Text="522×4="},
    @{Row=5;  Col=3; Text="690×3="},
    @{Row=5;  Col=4; Text="353×2="},
    @{Row=5;  Col=5; Text="160×5="},

    @{Row=10; Col=1; Text="645×4="},
    @{Row=10; Col=2; Text="906×7="},
    @{Row=10; Col=3; Text="531×2="},
    @{Row=10; Col=4; Text="487×4="},
    @{Row=10; Col=5; Text="726×5="},

    @{Row=15; Col=1; Text="318×5="},
    @{Row=15; Col=2; Text="869×4="},
    @{Row=15; Col=3; Text="216×9="},
    @{Row=15; Col=4; Text="891×8="},
    @{Row=15; Col=5; Text="376×3="},

    @{Row=20; Col=1; Text="517×2="},
    @{Row=20; Col=2; Text="358×2="},
    @{Row=20; Col=3; Text="537×3="},
    @{Row=20; Col=4; Text="263×2="},
    @{Row=20; Col=5; Text="690×7="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
